$wb = $excel.ActiveWorkbook

# New identifiers / timestamps coming from the localization "Generate Report for Handback" run.
$newMd1       = "890bf02a-7015-4d00-bb1a-6839176694b2.md"
$newMd1Path   = "e2e\890bf02a-7015-4d00-bb1a-6839176694b2.md"
$newMd2       = "ffffd62b1585-2dfd-4a30-b8fd-56a12646f84d.md"
$newMd2Path   = "e2e\ffffd62b1585-2dfd-4a30-b8fd-56a12646f84d.md"

$overviewDate = "2016-08-17 04:56:20"

$zhXlf        = "890bf02a-7015-4d00-bb1a-6839176694b2.616751ed80cc651b098b1a20eafef2c15a2d15b2.zh-cn.xlf"
$zhHoDate     = "2016-08-17 04:56:15"
$zhHbDate     = "2016-08-17 04:56:31"

$deXlf        = "890bf02a-7015-4d00-bb1a-6839176694b2.616751ed80cc651b098b1a20eafef2c15a2d15b2.de-de.xlf"
$deHbDate     = "2016-08-17 04:56:39"

# Original (unchanged) external hyperlink targets, keyed by sheet + old display text.
$ol1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf8f43aa38dccff6ba8c4a7dba3034ff3d4c11c8/e2e/06a1b3ac-3f75-4ab1-9ede-a483454965e9.md"
$ol2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf8f43aa38dccff6ba8c4a7dba3034ff3d4c11c8/e2e/e1366036-8bf9-40d2-a3d8-5229379de03f.md"
$zhcn1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8df44ec99b49693f2f4a22afc32004d9322f9474/e2e/06a1b3ac-3f75-4ab1-9ede-a483454965e9.md"
$zhcn2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8df44ec99b49693f2f4a22afc32004d9322f9474/e2e/e1366036-8bf9-40d2-a3d8-5229379de03f.md"
$dede1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4ee85a75d9bb84ec061c791ef5a429f39d0f33ae/e2e/06a1b3ac-3f75-4ab1-9ede-a483454965e9.md"
$dede2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4ee85a75d9bb84ec061c791ef5a429f39d0f33ae/e2e/e1366036-8bf9-40d2-a3d8-5229379de03f.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newMd1Path
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newMd2Path
$wsOverview.Range("G3").Value = $overviewDate

# Rebuild hyperlinks so the cached "display" text reflects the new file names
# while keeping the original (unchanged) target URLs.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ol1, $null, $null, $newMd1Path)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ol2, $null, $null, $newMd2Path)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $zhHoDate
$wsZh.Range("I2").Value = $newMd1
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhHbDate

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhHoDate
$wsZh.Range("I3").Value = $newMd2
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $zhHbDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $ol1, $null, $null, $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhcn1, $null, $null, $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ol2, $null, $null, $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhcn2, $null, $null, $newMd2)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $overviewDate
$wsDe.Range("I2").Value = $newMd1
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deHbDate

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $overviewDate
$wsDe.Range("I3").Value = $newMd2
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $deHbDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $ol1, $null, $null, $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $dede1, $null, $null, $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ol2, $null, $null, $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $dede2, $null, $null, $newMd2)
